# Adds the "CUMPLIMIENTO MENSUAL" sheet (sheet3) to the workbook, reproducing
# the VENTAS POR GRUPO -> CUMPLIMIENTO MENSUAL report with PRESUPUESTO / VENTA /
# POR CUMPLIR / CUMPLIMIENTO columns for advisor LOZANO MOLINA TITO.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet goes after the last existing sheet ("VENTA MENSUAL").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# --- column widths -------------------------------------------------------
# Excel's COM ColumnWidth setter pads the stored value by 5/6 of a character
# (default Calibri 11 / MDW=7), so the values below are chosen such that the
# round-tripped OOXML <col width="..."/> comes out to exactly 20/22/22/11/22/18.
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668
$ws.Columns.Item(3).ColumnWidth = 21.166666666666668
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668

# --- page margins ---------------------------------------------------------
$ps = $ws.PageSetup
$ps.LeftMargin = 0.75 * 72
$ps.RightMargin = 0.75 * 72
$ps.TopMargin = 1 * 72
$ps.BottomMargin = 1 * 72
$ps.HeaderMargin = 0.5 * 72
$ps.FooterMargin = 0.5 * 72

# --- header row (row 1) ----------------------------------------------------
# Reuse the existing bold/bordered header style from "VENTAS POR GRUPO" (sheet1)
# by copying formats only, then fill in this sheet's own header text.
$ws1.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats

$headers = @("ASESOR", "GRUPO", "PRESUPUESTO", "VENTA", "POR CUMPLIR", "CUMPLIMIENTO")
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# --- data rows (rows 2-18) --------------------------------------------------
$advisor = "LOZANO MOLINA TITO"
$groups = @(
    @("240X120 PORCELANATO", 344.284604629486),
    @("240X80 PORCELANATO", 3120.1145),
    @("FREGADEROS DE COCINA", 250.631825420901),
    @("GRANITO", 238.32),
    @("GRIFERIAS", 106.82),
    @("INODOROS", 560),
    @("LAVABOS", 625),
    @("LED", 300),
    @("NO RESURTIBLES", 650.25),
    @("OTROS", 0),
    @("PANELES DECORATIVOS", 350),
    @("PANELES PU", 230),
    @("PANELES PVC", 483),
    @("PIEDRA SINTERIZADA", 1638),
    @("PORCELANATO", 13061.58),
    @("PUERTAS DE SEGURIDAD", 342),
    @("SAL SOLUBLE", 1200)
)

$row = 2
foreach ($g in $groups) {
    $groupName = $g[0]
    $presupuesto = $g[1]

    $ws.Cells.Item($row, 1).Value = $advisor
    $ws.Cells.Item($row, 2).Value = $groupName

    $cPresupuesto = $ws.Cells.Item($row, 3)
    $cPresupuesto.Value = $presupuesto
    $cPresupuesto.NumberFormat = '"$"#,##0.00'

    $cVenta = $ws.Cells.Item($row, 4)
    $cVenta.Value = 0
    $cVenta.NumberFormat = '"$"#,##0.00'

    $cPorCumplir = $ws.Cells.Item($row, 5)
    $cPorCumplir.Value = $presupuesto
    $cPorCumplir.NumberFormat = '"$"#,##0.00'

    $cCumplimiento = $ws.Cells.Item($row, 6)
    $cCumplimiento.Value = 0
    $cCumplimiento.NumberFormat = "0.00%"

    $row++
}

# --- total row (row 19) -----------------------------------------------------
$totalRow = 19
$cTotalLabel = $ws.Cells.Item($totalRow, 2)
$cTotalLabel.Value = "TOTAL"
$cTotalLabel.HorizontalAlignment = -4152   # xlRight

$totalPresupuesto = 23500.00093005039

$cTotalPresupuesto = $ws.Cells.Item($totalRow, 3)
$cTotalPresupuesto.Value = $totalPresupuesto
$cTotalPresupuesto.NumberFormat = '"$"#,##0.00'

$cTotalVenta = $ws.Cells.Item($totalRow, 4)
$cTotalVenta.Value = 0
$cTotalVenta.NumberFormat = '"$"#,##0.00'

$cTotalPorCumplir = $ws.Cells.Item($totalRow, 5)
$cTotalPorCumplir.Value = $totalPresupuesto
$cTotalPorCumplir.NumberFormat = '"$"#,##0.00'

$cTotalCumplimiento = $ws.Cells.Item($totalRow, 6)
$cTotalCumplimiento.Value = 0
$cTotalCumplimiento.NumberFormat = "0.00%"

Write-Output "CUMPLIMIENTO MENSUAL sheet created"
